$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2024-07-12 16:41:23"
$ws.Range("B2").Value = "TestAaron Isac"
$ws.Range("C2").Value = "Portugal"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "+351924676500"
